# Update cryptos list data (prices and volume %) as scraped on Wed Dec  6 20:34:18 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "230.51", "43.978.12") but are
# stored as plain text in the sheet. Force text format first so Excel does not silently
# convert them into floating point numbers (which would lose trailing zeros / formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.978.12"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.260.92"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "230.51"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "64.24"
$ws.Range("E7").Value = "  +4.81%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +6.32%  "
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("D11").Value = "57.15"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "27.00"
$ws.Range("E12").Value = "  +14.40%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "2.599.85"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "15.73"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  +5.38%  "
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "2.255.88"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "43.918.42"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +7.49%  "
$ws.Range("D21").Value = "73.60"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "251.42"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("D28").Value = "3.25"
$ws.Range("E28").Value = "  +22.61%  "
$ws.Range("D29").Value = "171.06"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("D34").Value = "0.0705"
$ws.Range("E34").Value = "  +7.19%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("D37").Value = "3.78"
$ws.Range("E37").Value = "  +5.58%  "
$ws.Range("D38").Value = "6.49"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "0.0259"
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "0.000222"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").Value = "0.0977"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "17.36"
$ws.Range("E44").Value = "  +4.74%  "
$ws.Range("D45").Value = "8.22"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "98.14"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "4.44"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").Value = "10.32"
$ws.Range("E49").Value = "  +7.72%  "
$ws.Range("E50").Value = "  +5.68%  "
$ws.Range("D51").Value = "1.442.96"
$ws.Range("E51").Value = "  -1.88%  "

# Restore the default (no explicit number format) style on the price column so the
# saved workbook matches the original formatting.
$ws.Range("D2:D51").Style = "Normal"
